$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 24-56 (columns A-T), reflecting the refreshed
# weekly "Higo" price-report extract: existing rows 24-54 are overwritten
# with the updated dataset, and two additional rows (55-56) are appended,
# growing the sheet from A1:T54 to A1:T56.
$data = @(
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 45036, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 70, 20000, 20000, 20000, "`$/bandeja 7 kilos", "Región Metropolitana", 2857, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 45036, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 50, 14000, 14000, 14000, "`$/bandeja 7 kilos", "Región Metropolitana", 2000, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44302, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 340, 12000, 13000, 12500, "`$/bandeja 7 kilos", "Provincia de Santiago", 1786, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44685, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 100, 15000, 15000, 15000, "`$/bandeja 7 kilos", "Región Metropolitana", 2143, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44685, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 70, 12000, 12000, 12000, "`$/bandeja 7 kilos", "Región Metropolitana", 1714, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44306, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 50, 12000, 12000, 12000, "`$/bandeja 7 kilos", "Región Metropolitana", 1714, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44306, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 40, 9000, 9000, 9000, "`$/bandeja 7 kilos", "Región Metropolitana", 1286, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 45033, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 50, 20000, 20000, 20000, "`$/bandeja 7 kilos", "Región Metropolitana", 2857, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44664, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 80, 14000, 14000, 14000, "`$/bandeja 7 kilos", "Región Metropolitana", 2000, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44664, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 50, 12000, 12000, 12000, "`$/bandeja 7 kilos", "Región Metropolitana", 1714, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44322, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 100, 11000, 11000, 11000, "`$/bandeja 7 kilos", "Región Metropolitana", 1571, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 45020, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 120, 20000, 20000, 20000, "`$/bandeja 7 kilos", "Región Metropolitana", 2857, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44641, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 100, 13000, 13000, 13000, "`$/bandeja 7 kilos", "Región Metropolitana", 1857, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44644, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 85, 14000, 14000, 14000, "`$/bandeja 7 kilos", "Región Metropolitana", 2000, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 45029, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 120, 20000, 20000, 20000, "`$/bandeja 7 kilos", "Región Metropolitana", 2857, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44344, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 50, 12000, 12000, 12000, "`$/bandeja 7 kilos", "Región Metropolitana", 1714, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44300, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 150, 12000, 13000, 12500, "`$/bandeja 7 kilos", "Provincia de Santiago", 1786, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 45027, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 100, 20000, 20000, 20000, "`$/bandeja 7 kilos", "Región Metropolitana", 2857, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 45027, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 80, 14000, 14000, 14000, "`$/bandeja 7 kilos", "Región Metropolitana", 2000, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44699, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 50, 12000, 12000, 12000, "`$/bandeja 7 kilos", "Región Metropolitana", 1714, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44312, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 50, 13000, 13000, 13000, "`$/bandeja 7 kilos", "Región Metropolitana", 1857, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44312, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 20, 11000, 11000, 11000, "`$/bandeja 7 kilos", "Región Metropolitana", 1571, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44679, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 150, 12000, 12000, 12000, "`$/bandeja 7 kilos", "Región Metropolitana", 1714, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 45034, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 100, 20000, 20000, 20000, "`$/bandeja 7 kilos", "Región Metropolitana", 2857, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 45034, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 70, 14000, 14000, 14000, "`$/bandeja 7 kilos", "Región Metropolitana", 2000, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44694, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 80, 15000, 15000, 15000, "`$/bandeja 7 kilos", "Región Metropolitana", 2143, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44694, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 75, 12000, 12000, 12000, "`$/bandeja 7 kilos", "Región Metropolitana", 1714, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44687, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 100, 15000, 15000, 15000, "`$/bandeja 7 kilos", "Región Metropolitana", 2143, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44687, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 75, 12000, 12000, 12000, "`$/bandeja 7 kilos", "Región Metropolitana", 1714, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44987, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 120, 18000, 18000, 18000, "`$/bandeja 7 kilos", "Provincia de Santiago", 2571, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44316, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Primera", 40, 13000, 13000, 13000, "`$/bandeja 7 kilos", "Región Metropolitana", 1857, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44316, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 50, 11000, 11000, 11000, "`$/bandeja 7 kilos", "Región Metropolitana", 1571, 7),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44342, 13, "Fruta", 100101, "Berries", 100101006, "Higo", "Sin especificar", "Segunda", 50, 12000, 12000, 12000, "`$/bandeja 7 kilos", "Región Metropolitana", 1714, 7)
)

$startRow = 24
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
